$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New input values entered for the first two rows (2 stages) ---
$ws.Range("B2").Value = 5785.6263429175469
$ws.Range("D2").Value = 27.509567499970444
$ws.Range("F2").Value = 0.14999999999999999
$ws.Range("H2").Value = 10
$ws.Range("J2").Value = 0.0022499999999999998

$ws.Range("B3").Value = 6000.8168203639143
$ws.Range("D3").Value = 24.049923240085576
$ws.Range("F3").Value = 0.14999999999999999
$ws.Range("H3").Value = 10
$ws.Range("J3").Value = 0.059999999999999998

# --- Narrow the input columns to fit the new compact values ---
$ws.Columns.Item(2).ColumnWidth = 10.877604166666666
$ws.Columns.Item(4).ColumnWidth = 9.877604166666666
$ws.Columns.Item(6).ColumnWidth = 3.8776041666666665
$ws.Columns.Item(8).ColumnWidth = 2.3229166666666665
$ws.Columns.Item(10).ColumnWidth = 6.877604166666667

$wb.Save()
